$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Optimal_length_upravene") holds its numeric-looking values as
# text (shared strings), matching the rest of the sheet. A plain .Value
# assignment of a numeric-looking string would be auto-coerced to a real
# number, so the range is temporarily switched to Text format, the new,
# more precise values are written, and the original "0.00" number format
# (style index 1, already used by the sheet) is restored afterwards.
$rng = $ws.Range("C2:C14")
$rng.NumberFormat = "@"

$ws.Range("C2").Value = "0.152404782147249"
$ws.Range("C3").Value = "0.0902153118932702"
$ws.Range("C4").Value = "0.0746250795644468"
$ws.Range("C5").Value = "0.0673973478706693"
$ws.Range("C6").Value = "0.215514842938516"
$ws.Range("C7").Value = "0.108857430667365"
$ws.Range("C8").Value = "0.103602622664277"
$ws.Range("C9").Value = "0.116559476471344"
$ws.Range("C10").Value = "0.0418760373005853"
$ws.Range("C11").Value = "0.103022339810354"
$ws.Range("C12").Value = "0.0479718928123542"
$ws.Range("C13").Value = "0.0737665405180426"
$ws.Range("C14").Value = "0.125401365622548"

$rng.NumberFormat = "0.00"

# Move the active selection to C14, matching the final cursor position
$ws.Range("C14").Select()
